$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column D (batsman),
# shifting batsman..sr from D..I to F..K
$ws.Range("D:E").Insert()

# New header cells for the inserted columns
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data for existing rows (2 and 3)
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"

$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Kings XI Punjab"

# Append new row 4 with full data.
# Force the numeric-looking columns (totalRuns, totalBalls, total4s,
# total6s, sr) to be stored as text, matching the rest of the sheet.
$ws.Range("G4:K4").NumberFormat = "@"
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 12 2020"
$ws.Range("C4").Value = "RCB won by 82 runs"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
# Reuse the exact same batsman text already present in the sheet (same
# value/bytes used in rows 2 and 3) so the new row stays consistent.
$ws.Range("F4").Value = $ws.Range("F2").Value()
$ws.Range("G4").Value = "7"
$ws.Range("H4").Value = "10"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "70.00"
